$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2-246)
# from serial 45188 to serial 45189, matching the diff.
$ws.Range("C2:C246").Value = 45189
